$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new contingency rows ("line7"/"line8") are being inserted right
# after "line6" (row 7), ahead of the "extr" rows. Shift the existing
# "extr1".."extr8" rows (old rows 8-15) down by two rows (to new rows
# 10-17), carrying their name/from_bus/to_bus/in_service cells with them,
# before overwriting the now-vacated rows 8-9 with the new "line" entries.
$ws.Range("B8:E15").Copy($ws.Range("B10:E17"))

# Index column (A) keeps counting 0..15 down the whole table; propagate
# the existing bold/bordered/centered cell style used throughout column A
# onto the two newly-extended rows so no new (unused) style gets minted.
$ws.Range("A2:A3").Copy($ws.Range("A16:A17"))
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15

# --- New rows 8/9: "line7" / "line8" ---
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $false

$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# --- Re-run contingency results for the "extr" rows, now shifted to
#     rows 10-17; names already carried over by the copy above, only the
#     from_bus/to_bus/in_service outputs change ---
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $true

$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $false

$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
